$wb = $excel.ActiveWorkbook

# Rename the "All measured" sheet to "sce_prots_measured"
$allMeasured = $wb.Worksheets.Item("All measured")
$allMeasured.Name = "sce_prots_measured"

# Make it the active/selected sheet (matches tabSelected moving from
# "sce_SIG_Proteins_Osm 1" to "sce_prots_measured", and activeTab moving
# from index 5 to index 6)
$allMeasured.Select()
$allMeasured.Activate()
